# "DB: 'Ruta' is created"
# Populate the sm_ruta sheet (100 INSERT rows) and make it the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sm_ruta")

# --- Data rows: A (ID_Ruta), B (ID_Item), C (Ruta / image URL) ---------
for ($i = 1; $i -le 100; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $i
    $ws.Cells.Item($r, 2).Value = $i
    $ws.Cells.Item($r, 3).Value = "http://localhost:3783/SM/img/item$i/"
}

# --- Column D: INSERT statement formula --------------------------------
# Reproduce the author's actual authoring sequence (type the formula in
# D2, then fill down in two passes) so the shared-formula grouping in the
# saved XML matches: D2 alone, then D3:D66 (si=0), then D67:D101 (si=1).
$ws.Range("D2").Formula = '=CONCATENATE("INSERT INTO sm_ruta VALUES (",A2,",",B2,",",,"""",C2,""");")'
$ws.Range("D3:D66").Formula = '=CONCATENATE("INSERT INTO sm_ruta VALUES (",A3,",",B3,",",,"""",C3,""");")'
$ws.Range("D67:D101").Formula = '=CONCATENATE("INSERT INTO sm_ruta VALUES (",A67,",",B67,",",,"""",C67,""");")'

# --- Make sm_ruta the active sheet with A2:D101 selected ----------------
$ws.Select()
$ws.Range("A2:D101").Select()
